$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 19 for the "necrot3" back item
$ws.Range("A19").Value = "it_eq_back_necrot3"
$ws.Range("B19").Value = "itd_back_necrot3"
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = "50, 200"
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 20
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 10
$ws.Range("M19").Value = 10
$ws.Range("O19").Value = 10
$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 10
$ws.Range("Y19").Value = "res/assets/equipment/back/spritesheet_back_necrot3.png"

# Update the view state to match the recorded selection after editing
$excel.ActiveWindow.ScrollColumn = 17
$ws.Range("AA23").Select()
